$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 751, shifting the existing rows 751-767 down to 754-770
$ws.Rows("751:753").Insert()

# New row 751: La Araucanía / Vega Modelo de Temuco - Coliflor, week of 45239, Provincia del Elquí
$ws.Cells.Item(751, 1).Value = 10
$ws.Cells.Item(751, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(751, 3).Value = "La Araucanía"
$ws.Cells.Item(751, 4).Value = 45239
$ws.Cells.Item(751, 5).Value = 9
$ws.Cells.Item(751, 6).Value = 100112008
$ws.Cells.Item(751, 7).Value = "Coliflor"
$ws.Cells.Item(751, 8).Value = "Sin especificar"
$ws.Cells.Item(751, 9).Value = "Primera"
$ws.Cells.Item(751, 10).Value = 1800
$ws.Cells.Item(751, 11).Value = 1300
$ws.Cells.Item(751, 12).Value = 1300
$ws.Cells.Item(751, 13).Value = 1300
$ws.Cells.Item(751, 14).Value = "$/unidad"
$ws.Cells.Item(751, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(751, 16).Value = 1300
$ws.Cells.Item(751, 17).Value = 1
$ws.Cells.Item(751, 18).Value = "Hortaliza"

# New row 752: same week, Región del Maule
$ws.Cells.Item(752, 1).Value = 10
$ws.Cells.Item(752, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(752, 3).Value = "La Araucanía"
$ws.Cells.Item(752, 4).Value = 45239
$ws.Cells.Item(752, 5).Value = 9
$ws.Cells.Item(752, 6).Value = 100112008
$ws.Cells.Item(752, 7).Value = "Coliflor"
$ws.Cells.Item(752, 8).Value = "Sin especificar"
$ws.Cells.Item(752, 9).Value = "Primera"
$ws.Cells.Item(752, 10).Value = 1500
$ws.Cells.Item(752, 11).Value = 1400
$ws.Cells.Item(752, 12).Value = 1400
$ws.Cells.Item(752, 13).Value = 1400
$ws.Cells.Item(752, 14).Value = "$/unidad"
$ws.Cells.Item(752, 15).Value = "Región del Maule"
$ws.Cells.Item(752, 16).Value = 1400
$ws.Cells.Item(752, 17).Value = 1
$ws.Cells.Item(752, 18).Value = "Hortaliza"

# New row 753: same week, Segunda quality, Región Metropolitana
$ws.Cells.Item(753, 1).Value = 10
$ws.Cells.Item(753, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(753, 3).Value = "La Araucanía"
$ws.Cells.Item(753, 4).Value = 45239
$ws.Cells.Item(753, 5).Value = 9
$ws.Cells.Item(753, 6).Value = 100112008
$ws.Cells.Item(753, 7).Value = "Coliflor"
$ws.Cells.Item(753, 8).Value = "Sin especificar"
$ws.Cells.Item(753, 9).Value = "Segunda"
$ws.Cells.Item(753, 10).Value = 800
$ws.Cells.Item(753, 11).Value = 900
$ws.Cells.Item(753, 12).Value = 1000
$ws.Cells.Item(753, 13).Value = 938
$ws.Cells.Item(753, 14).Value = "$/unidad"
$ws.Cells.Item(753, 15).Value = "Región Metropolitana"
$ws.Cells.Item(753, 16).Value = 938
$ws.Cells.Item(753, 17).Value = 1
$ws.Cells.Item(753, 18).Value = "Hortaliza"
